$d = $word.ActiveDocument

$replacements = @(
    @{old="320×2=640"; new="588×8=4704"},
    @{old="804×5=4020"; new="836×6=5016"},
    @{old="308×7=2156"; new="499×8=3992"},
    @{old="721×3=2163"; new="151×3=453"},
    @{old="526×4=2104"; new="590×5=2950"},
    @{old="948×9=8532"; new="401×9=3609"},
    @{old="829×3=2487"; new="124×8=992"},
    @{old="656×2=1312"; new="236×5=1180"},
    @{old="576×2=1152"; new="311×2=622"},
    @{old="214×7=1498"; new="584×7=4088"},
    @{old="299×5=1495"; new="486×3=1458"},
    @{old="424×4=1696"; new="479×6=2874"},
    @{old="314×5=1570"; new="359×5=1795"},
    @{old="662×6=3972"; new="607×6=3642"},
    @{old="823×8=6584"; new="817×2=1634"},
    @{old="423×9=3807"; new="201×6=1206"},
    @{old="538×6=3228"; new="248×6=1488"},
    @{old="748×2=1496"; new="802×9=7218"},
    @{old="863×5=4315"; new="794×3=2382"},
    @{old="109×6=654"; new="143×6=858"},
    @{old="592×7=4144"; new="252×2=504"},
    @{old="365×7=2555"; new="823×7=5761"},
    @{old="605×6=3630"; new="510×5=2550"},
    @{old="895×6=5370"; new="733×2=1466"},
    @{old="162×4=648"; new="393×6=2358"}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $r.new, 2)
}
